$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 865.25
$ws.Cells.Item(32, 9).Value = 798.375
$ws.Cells.Item(32, 11).Value = 798.375
$ws.Cells.Item(32, 13).Value = -472.375
$ws.Cells.Item(48, 8).Value = 2999.6667
$ws.Cells.Item(48, 10).Value = 2999.6667
$ws.Cells.Item(48, 12).Value = 8999.000100000001
$ws.Cells.Item(48, 14).Value = -9583.000100000001
$ws.Cells.Item(56, 8).Value = 2999.6667
$ws.Cells.Item(56, 10).Value = 2999.6667
$ws.Cells.Item(56, 12).Value = 8999.000100000001
$ws.Cells.Item(56, 14).Value = -10067.0001
$ws.Cells.Item(62, 8).Value = 3336.75
$ws.Cells.Item(62, 9).Value = 2898
$ws.Cells.Item(62, 10).Value = 3600
$ws.Cells.Item(62, 11).Value = 2898
$ws.Cells.Item(62, 12).Value = 3600
$ws.Cells.Item(62, 13).Value = -2274
$ws.Cells.Item(62, 14).Value = -4848
$ws.Cells.Item(65, 8).Value = 3336.75
$ws.Cells.Item(65, 9).Value = 2898
$ws.Cells.Item(65, 10).Value = 3600
$ws.Cells.Item(65, 11).Value = 14490
$ws.Cells.Item(65, 12).Value = 18000
$ws.Cells.Item(65, 13).Value = -11370
$ws.Cells.Item(65, 14).Value = -24240
$ws.Cells.Item(69, 8).Value = 0
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 14).Value = ""
$ws.Cells.Item(72, 8).Value = 0
$ws.Cells.Item(72, 10).Value = 0
$ws.Cells.Item(72, 12).Value = 0
$ws.Cells.Item(72, 14).Value = ""
$ws.Cells.Item(74, 8).Value = 3499
$ws.Cells.Item(74, 9).Value = 4000
$ws.Cells.Item(74, 11).Value = 4000
$ws.Cells.Item(74, 13).Value = -3064
$ws.Cells.Item(77, 8).Value = 3499
$ws.Cells.Item(77, 9).Value = 4000
$ws.Cells.Item(77, 11).Value = 20000
$ws.Cells.Item(77, 13).Value = -15320
$ws.Cells.Item(106, 8).Value = 10000
$ws.Cells.Item(106, 9).Value = 0
$ws.Cells.Item(106, 11).Value = 0
$ws.Cells.Item(106, 13).Value = ""
$ws.Cells.Item(137, 8).Value = 1616.75
$ws.Cells.Item(137, 10).Value = 2240
$ws.Cells.Item(137, 12).Value = 6720
$ws.Cells.Item(137, 14).Value = -11820
$ws.Cells.Item(138, 8).Value = 1431.3529
$ws.Cells.Item(138, 10).Value = 3000
$ws.Cells.Item(138, 12).Value = 9000
$ws.Cells.Item(138, 14).Value = -19280
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(24, 8).Value = 40354
$ws.Cells.Item(24, 10).Value = 40354
$ws.Cells.Item(24, 12).Value = 40354
$ws.Cells.Item(24, 14).Value = -41102
$ws.Cells.Item(32, 8).Value = 3440.9348
$ws.Cells.Item(32, 9).Value = 3387.6904
$ws.Cells.Item(32, 11).Value = 3387.6904
$ws.Cells.Item(32, 13).Value = -3100.6904
$ws.Cells.Item(33, 8).Value = 28556.555
$ws.Cells.Item(33, 9).Value = 19498.75
$ws.Cells.Item(33, 10).Value = 35802.8
$ws.Cells.Item(33, 11).Value = 19498.75
$ws.Cells.Item(33, 12).Value = 35802.8
$ws.Cells.Item(33, 13).Value = -19169.75
$ws.Cells.Item(33, 14).Value = -36460.8
$ws.Cells.Item(45, 8).Value = 1784.5
$ws.Cells.Item(45, 9).Value = 1705
$ws.Cells.Item(45, 10).Value = 2500
$ws.Cells.Item(45, 11).Value = 1705
$ws.Cells.Item(45, 12).Value = 2500
$ws.Cells.Item(45, 13).Value = -1328
$ws.Cells.Item(45, 14).Value = -3254
$ws.Cells.Item(61, 8).Value = 3941.1667
$ws.Cells.Item(61, 9).Value = 3662.25
$ws.Cells.Item(61, 11).Value = 3662.25
$ws.Cells.Item(61, 13).Value = -3450.25
$ws.Cells.Item(74, 8).Value = 8000
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 13).Value = ""
$ws.Cells.Item(77, 8).Value = 8000
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 13).Value = ""
$ws.Cells.Item(100, 8).Value = 40354
$ws.Cells.Item(100, 10).Value = 40354
$ws.Cells.Item(100, 12).Value = 40354
$ws.Cells.Item(100, 14).Value = -42518
$ws.Cells.Item(136, 8).Value = 3941.1667
$ws.Cells.Item(136, 9).Value = 3662.25
$ws.Cells.Item(136, 11).Value = 10986.75
$ws.Cells.Item(136, 13).Value = -8436.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 531.8333
$ws.Cells.Item(20, 9).Value = 554.4
$ws.Cells.Item(20, 10).Value = 419
$ws.Cells.Item(20, 11).Value = 554.4
$ws.Cells.Item(20, 12).Value = 419
$ws.Cells.Item(20, 13).Value = -307.4
$ws.Cells.Item(20, 14).Value = -913
$ws.Cells.Item(82, 8).Value = 26375.584
$ws.Cells.Item(85, 8).Value = 26375.584
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 86.40000000000001
$ws.Cells.Item(7, 10).Value = 45
$ws.Cells.Item(7, 12).Value = 45
$ws.Cells.Item(7, 14).Value = -271
$ws.Cells.Item(31, 8).Value = 3524.5186
$ws.Cells.Item(31, 9).Value = 2465.8096
$ws.Cells.Item(31, 11).Value = 2465.8096
$ws.Cells.Item(31, 13).Value = -2170.8096
$ws.Cells.Item(34, 8).Value = 3524.5186
$ws.Cells.Item(34, 9).Value = 2465.8096
$ws.Cells.Item(34, 11).Value = 2465.8096
$ws.Cells.Item(34, 13).Value = -2263.8096
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 1775.7037
$ws.Cells.Item(4, 9).Value = 1908.125
$ws.Cells.Item(4, 11).Value = 5724.375
$ws.Cells.Item(4, 13).Value = -5612.375
$ws.Cells.Item(37, 8).Value = 99997.5
$ws.Cells.Item(37, 10).Value = 99997.5
$ws.Cells.Item(37, 12).Value = 299992.5
$ws.Cells.Item(37, 14).Value = -300216.5
$ws.Cells.Item(131, 8).Value = 1647
$ws.Cells.Item(131, 9).Value = 1163.5
$ws.Cells.Item(131, 11).Value = 3490.5
$ws.Cells.Item(131, 13).Value = 1549.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 11).Value = 0
$ws.Cells.Item(51, 13).Value = ""
$ws.Cells.Item(102, 8).Value = 759.8
$ws.Cells.Item(102, 9).Value = 759.8
$ws.Cells.Item(102, 11).Value = 759.8
$ws.Cells.Item(102, 13).Value = 862.2
$ws.Cells.Item(122, 8).Value = 9618265
$ws.Cells.Item(122, 9).Value = 13890520
$ws.Cells.Item(122, 10).Value = 5691.5
$ws.Cells.Item(122, 11).Value = 41671560
$ws.Cells.Item(122, 12).Value = 17074.5
$ws.Cells.Item(122, 13).Value = -41669110
$ws.Cells.Item(122, 14).Value = -21974.5
$ws.Cells.Item(132, 8).Value = 2661.75
$ws.Cells.Item(132, 9).Value = 2823.7
$ws.Cells.Item(132, 10).Value = 1852
$ws.Cells.Item(132, 11).Value = 8471.099999999999
$ws.Cells.Item(132, 12).Value = 5556
$ws.Cells.Item(132, 13).Value = -5941.099999999999
$ws.Cells.Item(132, 14).Value = -10616
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 3703.4
$ws.Cells.Item(122, 9).Value = 3503
$ws.Cells.Item(122, 11).Value = 10509
$ws.Cells.Item(122, 13).Value = -8059
$ws.Cells.Item(132, 8).Value = 9785.714
$ws.Cells.Item(132, 9).Value = 9833.333000000001
$ws.Cells.Item(132, 10).Value = 9750
$ws.Cells.Item(132, 11).Value = 29499.999
$ws.Cells.Item(132, 12).Value = 29250
$ws.Cells.Item(132, 13).Value = -26969.999
$ws.Cells.Item(132, 14).Value = -34310
$ws.Cells.Item(136, 8).Value = 2836
$ws.Cells.Item(136, 9).Value = 2836
$ws.Cells.Item(136, 11).Value = 8508
$ws.Cells.Item(136, 13).Value = -5958
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(58, 8).Value = 4666.6665
$ws.Cells.Item(58, 10).Value = 4000
$ws.Cells.Item(58, 12).Value = 4000
$ws.Cells.Item(58, 14).Value = -4616
$ws.Cells.Item(107, 8).Value = 207.88889
$ws.Cells.Item(107, 9).Value = 171.375
$ws.Cells.Item(107, 11).Value = 514.125
$ws.Cells.Item(107, 13).Value = 1405.875
$ws.Cells.Item(132, 8).Value = 2467
$ws.Cells.Item(132, 9).Value = 2467
$ws.Cells.Item(132, 11).Value = 7401
$ws.Cells.Item(132, 13).Value = -4871
$ws.Cells.Item(136, 8).Value = 13984.667
$ws.Cells.Item(136, 9).Value = 13987.5
$ws.Cells.Item(136, 10).Value = 13979
$ws.Cells.Item(136, 11).Value = 41962.5
$ws.Cells.Item(136, 12).Value = 41937
$ws.Cells.Item(136, 13).Value = -39412.5
$ws.Cells.Item(136, 14).Value = -47037
